$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 5174
$ws.Range("J2").Value = 5011
$ws.Range("L2").Value = 5011
$ws.Range("N2").Value = -5237

$ws.Range("H51").Value = 25854.125
$ws.Range("I51").Value = 10698.5
$ws.Range("J51").Value = 28019.215
$ws.Range("K51").Value = 10698.5
$ws.Range("L51").Value = 28019.215
$ws.Range("M51").Value = -10214.5
$ws.Range("N51").Value = -28987.215

$ws.Range("H132").Value = 1433.8823
$ws.Range("I132").Value = 917.0465
$ws.Range("J132").Value = 4211.875
$ws.Range("K132").Value = 2751.1395
$ws.Range("L132").Value = 12635.625
$ws.Range("M132").Value = -221.1395000000002
$ws.Range("N132").Value = -17695.625

$ws.Range("H138").Value = 9091.674999999999
$ws.Range("I138").Value = 5798.7
$ws.Range("J138").Value = 10089.546
$ws.Range("K138").Value = 17396.1
$ws.Range("L138").Value = 30268.638
$ws.Range("M138").Value = -12256.1
$ws.Range("N138").Value = -40548.638

$ws.Range("H141").Value = 2028.1666
$ws.Range("I141").Value = 1794.1904
$ws.Range("J141").Value = 3666
$ws.Range("K141").Value = 5382.5712
$ws.Range("L141").Value = 10998
$ws.Range("M141").Value = -202.5712000000003
$ws.Range("N141").Value = -21358

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 30797.578
$ws.Range("I32").Value = 29724.83
$ws.Range("K32").Value = 29724.83
$ws.Range("M32").Value = -29437.83

$ws.Range("H61").Value = 10175.944
$ws.Range("I61").Value = 4759.7
$ws.Range("J61").Value = 16946.25
$ws.Range("K61").Value = 4759.7
$ws.Range("L61").Value = 16946.25
$ws.Range("M61").Value = -4547.7
$ws.Range("N61").Value = -17370.25

$ws.Range("H74").Value = 387210.53
$ws.Range("I74").Value = 528451.8
$ws.Range("J74").Value = 3841.4285
$ws.Range("K74").Value = 528451.8
$ws.Range("L74").Value = 3841.4285
$ws.Range("M74").Value = -527577.8
$ws.Range("N74").Value = -5589.4285

$ws.Range("H77").Value = 387210.53
$ws.Range("I77").Value = 528451.8
$ws.Range("J77").Value = 3841.4285
$ws.Range("K77").Value = 2642259
$ws.Range("L77").Value = 19207.1425
$ws.Range("M77").Value = -2637891
$ws.Range("N77").Value = -27943.1425

$ws.Range("H110").Value = 6759181
$ws.Range("I110").Value = 9617142
$ws.Range("K110").Value = 9617142
$ws.Range("M110").Value = -9615097

$ws.Range("H122").Value = 3113.3225
$ws.Range("I122").Value = 2338.6875
$ws.Range("K122").Value = 7016.0625
$ws.Range("M122").Value = -4566.0625

$ws.Range("H125").Value = 60905
$ws.Range("J125").Value = 60905
$ws.Range("L125").Value = 60905
$ws.Range("N125").Value = -70745

$ws.Range("H132").Value = 10662.763
$ws.Range("I132").Value = 3058.963
$ws.Range("K132").Value = 9176.889000000001
$ws.Range("M132").Value = -6646.889000000001

$ws.Range("H136").Value = 10175.944
$ws.Range("I136").Value = 4759.7
$ws.Range("J136").Value = 16946.25
$ws.Range("K136").Value = 14279.1
$ws.Range("L136").Value = 50838.75
$ws.Range("M136").Value = -11729.1
$ws.Range("N136").Value = -55938.75

$ws.Range("H139").Value = 97571.336
$ws.Range("J139").Value = 97571.336
$ws.Range("L139").Value = 97571.336
$ws.Range("N139").Value = -107851.336

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H33").Value = 5699.75
$ws.Range("I33").Value = 5699.75
$ws.Range("K33").Value = 5699.75
$ws.Range("M33").Value = -5363.75

$ws.Range("H57").Value = 99999
$ws.Range("J57").Value = 99999
$ws.Range("L57").Value = 99999
$ws.Range("N57").Value = -101439

$ws.Range("H74").Value = 93950
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 93950
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 93950
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -95822

$ws.Range("H77").Value = 93950
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 93950
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 281850
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -291210

$ws.Range("H99").Value = 17565.562
$ws.Range("I99").Value = 22545.834
$ws.Range("K99").Value = 22545.834
$ws.Range("M99").Value = -21047.834

$ws.Range("H102").Value = 11836
$ws.Range("I102").Value = 11836
$ws.Range("K102").Value = 11836
$ws.Range("M102").Value = -8591

$ws.Range("H136").Value = 99999
$ws.Range("J136").Value = 99999
$ws.Range("L136").Value = 99999
$ws.Range("N136").Value = -110199

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9340.214
$ws.Range("J31").Value = 9459.462
$ws.Range("L31").Value = 9459.462
$ws.Range("N31").Value = -10049.462

$ws.Range("H34").Value = 9340.214
$ws.Range("J34").Value = 9459.462
$ws.Range("L34").Value = 9459.462
$ws.Range("N34").Value = -9863.462

$ws.Range("H122").Value = 1271.625
$ws.Range("I122").Value = 1167.5714
$ws.Range("K122").Value = 3502.7142
$ws.Range("M122").Value = -1052.7142

$ws.Range("H141").Value = 431240.88
$ws.Range("J141").Value = 448630.94
$ws.Range("L141").Value = 448630.94
$ws.Range("N141").Value = -458990.94

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 199.75
$ws.Range("I33").Value = 35.142857
$ws.Range("J33").Value = 327.77777
$ws.Range("K33").Value = 210.857142
$ws.Range("L33").Value = 1966.66662
$ws.Range("M33").Value = 72.14285799999999
$ws.Range("N33").Value = -2532.66662

$ws.Range("H107").Value = 784.8182
$ws.Range("J107").Value = 1320.1111
$ws.Range("L107").Value = 3960.3333
$ws.Range("N107").Value = -7800.3333

$ws.Range("H121").Value = 15874740
$ws.Range("J121").Value = 27779364
$ws.Range("L121").Value = 83338092
$ws.Range("N121").Value = -83340712

$ws.Range("H128").Value = 339496.72
$ws.Range("I128").Value = 339496.72
$ws.Range("K128").Value = 1018490.16
$ws.Range("M128").Value = -1013510.16

$ws.Range("H131").Value = 27800564
$ws.Range("I131").Value = 111112110
$ws.Range("J131").Value = 30046.666
$ws.Range("K131").Value = 333336330
$ws.Range("L131").Value = 90139.99800000001
$ws.Range("M131").Value = -333331290
$ws.Range("N131").Value = -100219.998

$ws.Range("H137").Value = 17571.428
$ws.Range("J137").Value = 3833.5
$ws.Range("L137").Value = 11500.5
$ws.Range("N137").Value = -21700.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H99").Value = 16556.666
$ws.Range("I99").Value = 16556.666
$ws.Range("K99").Value = 16556.666
$ws.Range("M99").Value = -14310.666

$ws.Range("H102").Value = 2021.1621
$ws.Range("I102").Value = 1864.1666
$ws.Range("K102").Value = 1864.1666
$ws.Range("M102").Value = -242.1666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 12502851
$ws.Range("I100").Value = 25002580
$ws.Range("K100").Value = 25002580
$ws.Range("M100").Value = -25002039

$ws.Range("H122").Value = 21743762
$ws.Range("I122").Value = 27782542
$ws.Range("J122").Value = 4148.4
$ws.Range("K122").Value = 83347626
$ws.Range("L122").Value = 12445.2
$ws.Range("M122").Value = -83345176
$ws.Range("N122").Value = -17345.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 9849.743
$ws.Range("I81").Value = 6797.125
$ws.Range("J81").Value = 12420.368
$ws.Range("K81").Value = 13594.25
$ws.Range("L81").Value = 24840.736
$ws.Range("M81").Value = -12533.25
$ws.Range("N81").Value = -26962.736

$ws.Range("H84").Value = 9849.743
$ws.Range("I84").Value = 6797.125
$ws.Range("J84").Value = 12420.368
$ws.Range("K84").Value = 67971.25
$ws.Range("L84").Value = 124203.68
$ws.Range("M84").Value = -62667.25
$ws.Range("N84").Value = -134811.68

$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()

$ws.Range("H132").Value = 5587.1523
$ws.Range("I132").Value = 4395.9653
$ws.Range("J132").Value = 7619.1763
$ws.Range("K132").Value = 13187.8959
$ws.Range("L132").Value = 22857.5289
$ws.Range("M132").Value = -10657.8959
$ws.Range("N132").Value = -27917.5289

$ws.Range("H136").Value = 3683.5173
$ws.Range("I136").Value = 2872.238
$ws.Range("K136").Value = 8616.714
$ws.Range("M136").Value = -6066.714
